$d = $word.ActiveDocument
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Project outcomes</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Two folders</w:t></w:r><w:r><w:t xml:space="preserve"> of generated slice images can be found in the attached zip. Cross- sectional images can be generated for any point of the skull, which is something that has never been done automatically before. </w:t></w:r></w:p><w:p><w:r><w:t>At the current stage of the project. A human has analyzed the generated images, and has produce</w:t></w:r><w:r><w:t>d</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>an age estimation, more accurate than estimates produced by already existing studies.</w:t></w:r></w:p><w:p><w:r><w:t>The automatically executed metric system is currently being worked on by me, and has the potential to introduce a significantly more efficient and effectiv</w:t></w:r><w:r><w:t>e</w:t></w:r><w:r><w:t xml:space="preserve"> than any existing</w:t></w:r><w:r><w:t xml:space="preserve"> one</w:t></w:r><w:r><w:t xml:space="preserve"> method for age estimation based on cranial suture</w:t></w:r><w:r><w:t xml:space="preserve"> analysis.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The project has so far produced more than 6500 images of cross-sectional suture images from more than 15 skulls. </w:t></w:r></w:p><w:p><w:r><w:t>After appropriate labeling, t</w:t></w:r><w:r><w:t xml:space="preserve">hose images will be used to train an image segmentation neural network, </w:t></w:r><w:r><w:t xml:space="preserve">to </w:t></w:r><w:r><w:t xml:space="preserve"> distinguish and </w:t></w:r><w:r><w:t xml:space="preserve">mark </w:t></w:r><w:r><w:t xml:space="preserve">a </w:t></w:r><w:r><w:t>region in which</w:t></w:r><w:r><w:t>, a</w:t></w:r><w:r><w:t xml:space="preserve"> part of a suture is present</w:t></w:r><w:r><w:t>. The</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>properties</w:t></w:r><w:r><w:t xml:space="preserve"> of the</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>pointed-out</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>suture</w:t></w:r><w:r><w:t xml:space="preserve"> (suture color, suture width) will further be analyzed. An overall statistic and conclusion about the age of the individual will be made based on the results produced </w:t></w:r><w:r><w:t>on</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> each image.</w:t></w:r></w:p><w:p><w:r><w:t>An example of what the image segmentation neural network has to achieve can be found in</w:t></w:r><w:r><w:t xml:space="preserve"> the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>project_dev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> folder.</w:t></w:r></w:p><w:p/>
'@
$d.Content.InsertXML($xml)
